# Auto-generated edit script: updates numeric cells in Sheets/Anima_Profits.xlsx
# per the commit 'chore: update Sheets via scheduled runner'.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 152.35715
$ws.Range("I33").Value = 141
$ws.Range("K33").Value = 141
$ws.Range("M33").Value = 88
$ws.Range("H116").Value = 2434.1667
$ws.Range("I116").Value = 2281
$ws.Range("J116").Value = 3200
$ws.Range("K116").Value = 2281
$ws.Range("L116").Value = 3200
$ws.Range("M116").Value = 1161
$ws.Range("N116").Value = -10084
$ws.Range("H137").Value = 3455.2903
$ws.Range("I137").Value = 3976.9375
$ws.Range("K137").Value = 11930.8125
$ws.Range("M137").Value = -9380.8125
$ws.Range("H138").Value = 2036.9688
$ws.Range("I138").Value = 1405.5366
$ws.Range("J138").Value = 3162.5652
$ws.Range("K138").Value = 4216.6098
$ws.Range("L138").Value = 9487.695599999999
$ws.Range("M138").Value = 923.3901999999998
$ws.Range("N138").Value = -19767.6956
$ws.Range("H141").Value = 3542.1924
$ws.Range("I141").Value = 1609.2174
$ws.Range("J141").Value = 18361.666
$ws.Range("K141").Value = 4827.6522
$ws.Range("L141").Value = 55084.99800000001
$ws.Range("M141").Value = 352.3477999999996
$ws.Range("N141").Value = -65444.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 17750
$ws.Range("J22").Value = 17750
$ws.Range("L22").Value = 17750
$ws.Range("N22").Value = -18348
$ws.Range("H32").Value = 540597.1
$ws.Range("I32").Value = 605533.6
$ws.Range("J32").Value = 21105.25
$ws.Range("K32").Value = 605533.6
$ws.Range("L32").Value = 21105.25
$ws.Range("M32").Value = -605246.6
$ws.Range("N32").Value = -21679.25
$ws.Range("H38").Value = 3373
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H44").Value = 27005.375
$ws.Range("J44").Value = 27005.375
$ws.Range("L44").Value = 27005.375
$ws.Range("N44").Value = -27981.375
$ws.Range("H61").Value = 2289.8
$ws.Range("I61").Value = 1914.7142
$ws.Range("J61").Value = 3165
$ws.Range("K61").Value = 1914.7142
$ws.Range("L61").Value = 3165
$ws.Range("M61").Value = -1702.7142
$ws.Range("N61").Value = -3589
$ws.Range("H80").Value = 52499
$ws.Range("J80").Value = 52499
$ws.Range("L80").Value = 52499
$ws.Range("N80").Value = -54495
$ws.Range("H83").Value = 52499
$ws.Range("J83").Value = 52499
$ws.Range("L83").Value = 157497
$ws.Range("N83").Value = -167481
$ws.Range("H132").Value = 3522.7112
$ws.Range("I132").Value = 2442.258
$ws.Range("J132").Value = 5915.143
$ws.Range("K132").Value = 7326.773999999999
$ws.Range("L132").Value = 17745.429
$ws.Range("M132").Value = -4796.773999999999
$ws.Range("N132").Value = -22805.429
$ws.Range("H136").Value = 2289.8
$ws.Range("I136").Value = 1914.7142
$ws.Range("J136").Value = 3165
$ws.Range("K136").Value = 5744.142599999999
$ws.Range("L136").Value = 9495
$ws.Range("M136").Value = -3194.142599999999
$ws.Range("N136").Value = -14595

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1332.0625
$ws.Range("J80").Value = 227
$ws.Range("L80").Value = 227
$ws.Range("N80").Value = -2223
$ws.Range("H83").Value = 1332.0625
$ws.Range("J83").Value = 227
$ws.Range("L83").Value = 1135
$ws.Range("N83").Value = -11119
$ws.Range("H105").Value = 62501404
$ws.Range("I105").Value = 62501404
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 62501404
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 991
$ws.Range("I107").Value = 606.8889
$ws.Range("J107").Value = 1682.4
$ws.Range("K107").Value = 606.8889
$ws.Range("L107").Value = 1682.4
$ws.Range("M107").Value = 1313.1111
$ws.Range("N107").Value = -5522.4
$ws.Range("H134").Value = 3507.8845
$ws.Range("I134").Value = 2949.6875
$ws.Range("K134").Value = 8849.0625
$ws.Range("M134").Value = -6314.0625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 1011
$ws.Range("I25").Value = 1011
$ws.Range("K25").Value = 1011
$ws.Range("M25").Value = -837
$ws.Range("H31").Value = 7315.023
$ws.Range("I31").Value = 1571.9286
$ws.Range("K31").Value = 1571.9286
$ws.Range("M31").Value = -1276.9286
$ws.Range("H34").Value = 7315.023
$ws.Range("I34").Value = 1571.9286
$ws.Range("K34").Value = 1571.9286
$ws.Range("M34").Value = -1369.9286
$ws.Range("H57").Value = 18000
$ws.Range("J57").Value = 18000
$ws.Range("L57").Value = 18000
$ws.Range("N57").Value = -19120
$ws.Range("H107").Value = 2500558.5
$ws.Range("I107").Value = 6944721
$ws.Range("J107").Value = 717
$ws.Range("K107").Value = 6944721
$ws.Range("L107").Value = 717
$ws.Range("M107").Value = -6942801
$ws.Range("N107").Value = -4557
$ws.Range("H132").Value = 35355180
$ws.Range("I132").Value = 40001370
$ws.Range("J132").Value = 20835846
$ws.Range("K132").Value = 120004110
$ws.Range("L132").Value = 62507538
$ws.Range("M132").Value = -120001580
$ws.Range("N132").Value = -62512598
$ws.Range("H134").Value = 4311.8423
$ws.Range("I134").Value = 4024.4412
$ws.Range("J134").Value = 6754.75
$ws.Range("K134").Value = 12073.3236
$ws.Range("L134").Value = 20264.25
$ws.Range("M134").Value = -9538.3236
$ws.Range("N134").Value = -25334.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2833.3333
$ws.Range("I19").Value = 2166.6667
$ws.Range("J19").Value = 3500
$ws.Range("K19").Value = 6500.000100000001
$ws.Range("L19").Value = 10500
$ws.Range("M19").Value = -6326.000100000001
$ws.Range("N19").Value = -10848
$ws.Range("H131").Value = 1665.0834
$ws.Range("J131").Value = 1888.1
$ws.Range("L131").Value = 5664.299999999999
$ws.Range("N131").Value = -15744.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H132").Value = 3826.3125
$ws.Range("I132").Value = 3564.0417
$ws.Range("J132").Value = 4613.125
$ws.Range("K132").Value = 10692.1251
$ws.Range("L132").Value = 13839.375
$ws.Range("M132").Value = -8162.125100000001
$ws.Range("N132").Value = -18899.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2830
$ws.Range("H122").Value = 3295.6667
$ws.Range("I122").Value = 2750.4443
$ws.Range("J122").Value = 3704.5833
$ws.Range("K122").Value = 8251.332900000001
$ws.Range("L122").Value = 11113.7499
$ws.Range("M122").Value = -5801.332900000001
$ws.Range("N122").Value = -16013.7499
$ws.Range("H132").Value = 2535.5312
$ws.Range("I132").Value = 2116.2593
$ws.Range("K132").Value = 6348.777900000001
$ws.Range("M132").Value = -3818.777900000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9807886
$ws.Range("I132").Value = 5602.875
$ws.Range("J132").Value = 18521026
$ws.Range("K132").Value = 16808.625
$ws.Range("L132").Value = 55563078
$ws.Range("M132").Value = -14278.625
$ws.Range("N132").Value = -55568138
$ws.Range("H136").Value = 2337.6182
$ws.Range("I136").Value = 1898.5814
$ws.Range("K136").Value = 5695.7442
$ws.Range("M136").Value = -3145.7442
